$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.142.84"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "1.668.21"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("D5").Value = "209.37"
$ws.Range("E5").Value = "  -3.34%  "
$ws.Range("D6").Value = "0.5229"
$ws.Range("E6").Value = "  -0.82%  "
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("D8").Value = "0.2622"
$ws.Range("E8").Value = "  -2.87%  "
$ws.Range("D9").Value = "0.06328"
$ws.Range("E9").Value = "  -0.84%  "
$ws.Range("D10").Value = "21.18"
$ws.Range("E10").Value = "  -1.84%  "
$ws.Range("D11").Value = "0.07536"
$ws.Range("E11").Value = "  -1.60%  "
$ws.Range("D12").Value = "1.683.95"
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("D13").Value = "4.444"
$ws.Range("E13").Value = "  -1.47%  "
$ws.Range("D14").Value = "0.5497"
$ws.Range("E14").Value = "  -4.59%  "
$ws.Range("D15").Value = "66.35"
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("D16").Value = "0.000007954"
$ws.Range("D17").Value = "26.153.51"
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("D18").Value = "1.002"
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("D19").Value = "4.750"
$ws.Range("E19").Value = "  -2.52%  "
$ws.Range("D20").Value = "186.62"
$ws.Range("E20").Value = "  -1.57%  "
$ws.Range("E21").Value = "  -4.73%  "
$ws.Range("D22").Value = "6.185"
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("E23").Value = "  -0.51%  "
$ws.Range("D24").Value = "149.19"
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("D25").Value = "0.1248"
$ws.Range("E25").Value = "  -0.95%  "
$ws.Range("D26").Value = "7.500"
$ws.Range("E26").Value = "  -3.86%  "
$ws.Range("D27").Value = "15.86"
$ws.Range("E27").Value = "  +0.77%  "
$ws.Range("D28").Value = "0.06376"
$ws.Range("E28").Value = "  +1.71%  "
$ws.Range("D29").Value = "1.348"
$ws.Range("E29").Value = "  -1.65%  "
$ws.Range("D30").Value = "1.275"
$ws.Range("E30").Value = "  -3.39%  "
$ws.Range("D31").Value = "3.513"
$ws.Range("E31").Value = "  -1.66%  "
$ws.Range("D32").Value = "3.410"
$ws.Range("E32").Value = "  -4.28%  "
$ws.Range("D33").Value = "1.647"
$ws.Range("E33").Value = "  -2.36%  "
$ws.Range("D34").Value = "1.004"
$ws.Range("E34").Value = "  -1.94%  "
$ws.Range("E35").Value = "  -1.70%  "
$ws.Range("D36").Value = "2.407"
$ws.Range("E36").Value = "  -0.58%  "
$ws.Range("D37").Value = "2.738"
$ws.Range("E37").Value = "  -0.76%  "
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").Value = "6.147"
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "1.112.47"
$ws.Range("E39").Value = "  +1.10%  "
$ws.Range("E40").Value = "  -0.32%  "
$ws.Range("D41").Value = "0.8663"
$ws.Range("E41").Value = "  -3.53%  "
$ws.Range("D42").Value = "1.003"
$ws.Range("E42").Value = "  -0.76%  "
$ws.Range("D43").Value = "100.32"
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("D44").Value = "1.823.39"
$ws.Range("E44").Value = "  -0.57%  "
$ws.Range("D45").Value = "0.00000000111"
$ws.Range("E45").Value = "  +0.20%  "
$ws.Range("D46").Value = "55.49"
$ws.Range("E46").Value = "  -3.36%  "
$ws.Range("D47").Value = "1.004"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("D48").Value = "8.053"
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("D49").Value = "0.05234"
$ws.Range("E49").Value = "  -0.81%  "
$ws.Range("D50").Value = "0.4237"
$ws.Range("E50").Value = "  -1.18%  "
$ws.Range("D51").Value = "5.934"
$ws.Range("E51").Value = "  -1.44%  "
